# The deck originally carries two embedded DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (unused placeholder, only wired to the Notes Master)
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" (the theme actually driving the slide master)
# The target edit swaps the two themes' content so the deck's live design becomes the
# plain default "Office Theme" colors instead of "Integral" / "Red Violet".
#
# PowerPoint's object model edits theme colours through ThemeColorScheme (12 slots, in
# clrScheme document order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). Writing through
# any slide's ThemeColorScheme updates the theme bound to the presentation's slide master.

function Convert-RGB {
    param([int]$r, [int]$g, [int]$b)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette: the stock "Office" colour scheme (dk1..folHlink), which is the palette
# that used to live in theme1.xml and should now become the deck's active theme colours.
$officeColors = @(
    @(0x00, 0x00, 0x00),  # dk1
    @(0xFF, 0xFF, 0xFF),  # lt1
    @(0x44, 0x54, 0x6A),  # dk2
    @(0xE7, 0xE6, 0xE6),  # lt2
    @(0x5B, 0x9B, 0xD5),  # accent1
    @(0xED, 0x7D, 0x31),  # accent2
    @(0xA5, 0xA5, 0xA5),  # accent3
    @(0xFF, 0xC0, 0x00),  # accent4
    @(0x44, 0x72, 0xC4),  # accent5
    @(0x70, 0xAD, 0x47),  # accent6
    @(0x05, 0x63, 0xC1),  # hlink
    @(0x95, 0x4F, 0x72)   # folHlink
)

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $rgb = $officeColors[$i]
    $tcs.Item($i + 1).RGB = Convert-RGB $rgb[0] $rgb[1] $rgb[2]
}
